$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/19/2025  Through  5/25/2025"

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("N15").Value = -76.470588235294
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = 29.411764705882
$ws.Range("L16").Value = -10.204081632653
$ws.Range("M16").Value = -60.360360360360
$ws.Range("N16").Value = -93.519882179676
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 16.666666666666
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 28
$ws.Range("I17").Value = 144
$ws.Range("J17").Value = 104
$ws.Range("K17").Value = 38.461538461538
$ws.Range("L17").Value = 13.385826771653
$ws.Range("M17").Value = 35.849056603773
$ws.Range("N17").Value = -60.220994475138
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 28
$ws.Range("K18").Value = -7.142857142857
$ws.Range("L18").Value = -3.703703703703
$ws.Range("M18").Value = -78.151260504201
$ws.Range("N18").Value = -96.605744125326
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = 113.333333333333
$ws.Range("I19").Value = 124
$ws.Range("J19").Value = 101
$ws.Range("K19").Value = 22.772277227722
$ws.Range("L19").Value = -32.608695652173
$ws.Range("M19").Value = -18.421052631578
$ws.Range("N19").Value = -47.899159663865
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = 3.571428571428
$ws.Range("L20").Value = -29.268292682926
$ws.Range("M20").Value = -50.847457627118
$ws.Range("N20").Value = -94.649446494464
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 18.75
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 61
$ws.Range("H21").Value = 34.426229508196
$ws.Range("I21").Value = 376
$ws.Range("J21").Value = 310
$ws.Range("K21").Value = 21.290322580645
$ws.Range("L21").Value = -14.350797266514
$ws.Range("M21").Value = -32.007233273056
$ws.Range("N21").Value = -85.73055028463
$ws.Range("M22").Value = -92.857142857142
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 21.428571428571
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 35
$ws.Range("I24").Value = 380
$ws.Range("J24").Value = 363
$ws.Range("K24").Value = 4.683195592286
$ws.Range("L24").Value = -14.221218961625
$ws.Range("M24").Value = 2.702702702702
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 66.666666666666
$ws.Range("I25").Value = 110
$ws.Range("J25").Value = 68
$ws.Range("K25").Value = 61.764705882352
$ws.Range("L25").Value = 23.595505617977
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = -42.105263157894
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 47
$ws.Range("H26").Value = -25.531914893617
$ws.Range("I26").Value = 182
$ws.Range("J26").Value = 211
$ws.Range("K26").Value = -13.744075829383
$ws.Range("L26").Value = -10.344827586206
$ws.Range("M26").Value = -42.586750788643
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 33.333333333333
$ws.Range("J28").Value = 27
$ws.Range("K28").Value = -11.111111111111
$ws.Range("L28").Value = 41.176470588235
$ws.Range("N29").Value = -96.226415094339
$ws.Range("N30").Value = -96.153846153846

# --- Numeric -> Text ("0" / "***.*") conversions (style 13) ---
# Set the (possibly numeric-looking) text first (apostrophe forces text type),
# then copy number-format only from a stable style-13 donor cell (C14) so the
# cell reuses the existing style index instead of Excel minting a new one.
$ws.Range("C20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D31").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E31").PasteSpecial(-4122)

# --- Text -> Numeric conversions (style 14 = integer counts, style 15 = percents) ---
$ws.Range("D18").Value = 3
$ws.Range("C18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("K14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = -66.666666666666
$ws.Range("D28").Value = 1
$ws.Range("C18").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100

$excel.CutCopyMode = $false